# Fix the timing and typo in PauliTwoDesign ansatz
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Optimizer_Ansatz")
$ws.Activate()

# Widen column H to fit the new numbers (target OOXML width="20";
# this engine's ColumnWidth has a fixed +0.8333.. padding offset vs. the
# stored XML width, so back it out to land exactly on 20)
$ws.Columns.Item(8).ColumnWidth = 19.166666666666668

# Row 2 (EfficientSU2 / SPSA)
$ws.Range("F2").Value = 250.5
$ws.Range("G2").Value = 1.28
$ws.Range("H2").Value = 9.862
$ws.Range("I2").Value = 0.000826
$ws.Range("J2").Value = 9.862
$ws.Range("K2").Value = 0.00335
$ws.Range("L2").Value = "928KB"

# Row 3 (EfficientSU2 / COBYLA)
$ws.Range("F3").Value = 127.9
$ws.Range("G3").Value = 4.78
$ws.Range("L3").Value = "936KB"

# Row 4 (EfficientSU2 / SLSQP)
$ws.Range("F4").Value = 16.15
$ws.Range("G4").Value = 21.12
$ws.Range("L4").Value = "928KB"

# Row 5 (EfficientSU2 / P_BFGS)
$ws.Range("F5").Value = 33.37
$ws.Range("G5").Value = 282.84
$ws.Range("L5").Value = "8.56GB"

# Row 6 (EfficientSU2 / ADAM)
$ws.Range("F6").Value = 18.83
$ws.Range("G6").Value = 1778.83
$ws.Range("L6").Value = "241.67MB"

# Row 7 (TwoLocal / SPSA) - correct job id, fill in timing/memory
$ws.Range("A7").Value = 18994401
$ws.Range("F7").Value = 746.3
$ws.Range("G7").Value = 1.44
$ws.Range("L7").Value = "1.16MB"

# Row 8 (TwoLocal / COBYLA) - correct job id, fill in timing/memory
$ws.Range("A8").Value = 18994593
$ws.Range("F8").Value = 106.1
$ws.Range("G8").Value = 5.61
$ws.Range("L8").Value = "1.16MB"

# Row 9 (TwoLocal / SLSQP) - correct job id, fill in timing/memory
$ws.Range("A9").Value = 18994600
$ws.Range("F9").Value = 10.09
$ws.Range("G9").Value = 23.83
$ws.Range("L9").Value = "1.16MB"

# Row 10 (TwoLocal / P_BFGS) - add missing job id, fill in timing/memory
$ws.Range("A10").Value = 18994607
$ws.Range("F10").Value = 31.29
$ws.Range("G10").Value = 382.37
$ws.Range("L10").Value = "8.87GB"

# Row 11 (TwoLocal / ADAM) - add missing job id
$ws.Range("A11").Value = 18994609

# Update the view: zoom to 145% and move the selection to D7
$ws.Range("D7").Select()
$excel.ActiveWindow.Zoom = 145
